# Draft L2C Q4 2021 report
# Update the "L2C Group" summary table: refresh the n (percent) counts for
# each disposition category, insert a new (unlabeled) row for a second
# "Dropped" count of 2 (0.6), and refresh the Total row.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Update existing category counts (column 2) by row.
$t.Cell(2, 2).Range.Text = "101 (30.7)"   # UCM+SP
$t.Cell(3, 2).Range.Text = "102 (31.0)"   # UCM
$t.Cell(4, 2).Range.Text = "98 (29.8)"    # L2C
$t.Cell(5, 2).Range.Text = "23 (7.0)"     # NS V2
$t.Cell(6, 2).Range.Text = "2 (0.6)"      # Dropped

# Insert a new blank row right before the "Total" row (currently row 8)
# and fill in its second cell; the first cell is left empty.
$totalRow = $t.Rows.Item(8)
$newRow = $t.Rows.Add($totalRow)
$newRow.Cells.Item(2).Range.Text = "2 (0.6)"

# Refresh the Total row's count (now row 9 after the insertion).
$t.Cell($t.Rows.Count, 2).Range.Text = "329 (100.0)"
